$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.99999999672415563,
  0.99362820179469236,
  0.96501460430741304,
  0.9556732347197856,
  0.94675673838898322,
  0.92490765949875753,
  0.9243145106954912,
  0.92408456272513217,
  0.92775822365454053,
  0.93243306147388982,
  0.9332716165273085,
  0.93512662146314329,
  0.92384015411402165,
  0.919672870741161,
  0.91708140098769364,
  0.91457492705765353,
  0.91086699552647266,
  0.9097580675184338,
  0.99191021304280413,
  0.98479298406511395,
  0.98339446946454678,
  0.98212995941190107,
  0.9500942220036801,
  0.92813777273402986,
  0.9216803139429286,
  0.92125003782625381,
  0.91640082150818103,
  0.89491921226677706,
  0.87964447124441481,
  0.87307331476261418,
  0.86541880382244285,
  0.86373940940795824,
  0.86321938492077255
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
